$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SIQ")

$ws.Range("I11").Value = "TSH : after you agjust an press the mode you should save and back to watch mode "
$ws.Range("I12").Value = "TSH : when you enter alarm you just adjust the alarm then after finish you go back to normal watch then when alarm time comes you alert the user as normal watch do"
$ws.Range("I13").Value = "TSH: you should stop the alarm even after certain time if no button pressed could be 30 sec or press increment button "
$ws.Range("I14").Value = "TSH : don't get the question ?!"
$ws.Range("I15").Value = "TSH : the default is 12 and alaramed disabeld "
$ws.Range("I16").Value = "TSH : just like normal watch do , when increment the hours and you are 12 then will be 1 , nut if you can adjut the PM/AM then just increment the hour and don't togle the PM/AM as we adjust it seperately "
$ws.Range("I17").Value = "TSH : 30 sec "
